$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the label-column formatting (bold, thin border, center/top aligned)
# from the last existing styled row (J25) onto the six new rows being added
# (J26:J31), matching the style used by J3:J25.
$ws.Range("J25").Copy() | Out-Null
$ws.Range("J26:J31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"
$ws.Range("A3").Value = "crude"
$ws.Range("B3").Value = 0.8235294117647058
$ws.Range("C3").Value = 28
$ws.Range("D3").Value = 28
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 6
$ws.Range("J3").Value = "happy"
$ws.Range("K3").Value = 0.9615384615384616
$ws.Range("L3").Value = 25
$ws.Range("M3").Value = 25
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 1
$ws.Range("A4").Value = "fraud"
$ws.Range("B4").Value = 0.6944444444444444
$ws.Range("C4").Value = 25
$ws.Range("D4").Value = 25
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 11
$ws.Range("J4").Value = "love"
$ws.Range("K4").Value = 0.9347826086956522
$ws.Range("L4").Value = 43
$ws.Range("M4").Value = 43
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 3
$ws.Range("A5").Value = "crisis"
$ws.Range("B5").Value = 0.5958904109589042
$ws.Range("C5").Value = 174
$ws.Range("D5").Value = 174
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 118
$ws.Range("J5").Value = "best"
$ws.Range("K5").Value = 0.9322033898305084
$ws.Range("L5").Value = 55
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 4
$ws.Range("A6").Value = "panic"
$ws.Range("B6").Value = 0.187984496124031
$ws.Range("C6").Value = 97
$ws.Range("D6").Value = 97
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 419
$ws.Range("J6").Value = "interesting"
$ws.Range("K6").Value = 0.8787878787878788
$ws.Range("L6").Value = 29
$ws.Range("M6").Value = 29
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 4
$ws.Range("A7").Value = "sc"
$ws.Range("B7").Value = 0.1693121693121693
$ws.Range("C7").Value = 32
$ws.Range("D7").Value = 32
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 157
$ws.Range("J7").Value = "great"
$ws.Range("K7").Value = 0.8482142857142857
$ws.Range("L7").Value = 95
$ws.Range("M7").Value = 95
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 17
$ws.Range("J8").Value = "won"
$ws.Range("K8").Value = 0.8461538461538461
$ws.Range("L8").Value = 33
$ws.Range("M8").Value = 33
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 6
$ws.Range("J9").Value = "thanks"
$ws.Range("K9").Value = 0.7926829268292683
$ws.Range("L9").Value = 65
$ws.Range("M9").Value = 65
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 17
$ws.Range("J10").Value = "thank"
$ws.Range("K10").Value = 0.7890625
$ws.Range("L10").Value = 101
$ws.Range("M10").Value = 101
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 27
$ws.Range("J11").Value = "special"
$ws.Range("K11").Value = 0.7777777777777778
$ws.Range("L11").Value = 28
$ws.Range("M11").Value = 28
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 8
$ws.Range("J12").Value = "positive"
$ws.Range("K12").Value = 0.7758620689655172
$ws.Range("L12").Value = 45
$ws.Range("M12").Value = 45
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 13
$ws.Range("J13").Value = "safety"
$ws.Range("K13").Value = 0.7254901960784313
$ws.Range("L13").Value = 37
$ws.Range("M13").Value = 37
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 14
$ws.Range("J14").Value = "confidence"
$ws.Range("K14").Value = 0.7222222222222222
$ws.Range("L14").Value = 26
$ws.Range("M14").Value = 26
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 10
$ws.Range("J15").Value = "free"
$ws.Range("K15").Value = 0.7166666666666667
$ws.Range("L15").Value = 86
$ws.Range("M15").Value = 86
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 34
$ws.Range("J16").Value = "safe"
$ws.Range("K16").Value = 0.7112676056338029
$ws.Range("L16").Value = 101
$ws.Range("M16").Value = 101
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 41
$ws.Range("J17").Value = "good"
$ws.Range("K17").Value = 0.675
$ws.Range("L17").Value = 108
$ws.Range("M17").Value = 108
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 52
$ws.Range("J18").Value = "support"
$ws.Range("K18").Value = 0.660377358490566
$ws.Range("L18").Value = 70
$ws.Range("M18").Value = 70
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 36
$ws.Range("J19").Value = "fresh"
$ws.Range("K19").Value = 0.625
$ws.Range("L19").Value = 30
$ws.Range("M19").Value = 30
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 18
$ws.Range("J20").Value = "well"
$ws.Range("K20").Value = 0.6170212765957447
$ws.Range("L20").Value = 58
$ws.Range("M20").Value = 58
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 36
$ws.Range("J21").Value = "better"
$ws.Range("K21").Value = 0.6031746031746031
$ws.Range("L21").Value = 38
$ws.Range("M21").Value = 38
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 25
$ws.Range("J22").Value = "relief"
$ws.Range("K22").Value = 0.56
$ws.Range("L22").Value = 28
$ws.Range("M22").Value = 28
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 22
$ws.Range("J23").Value = "heroes"
$ws.Range("K23").Value = 0.5319148936170213
$ws.Range("L23").Value = 25
$ws.Range("M23").Value = 25
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 22
$ws.Range("J24").Value = "hand"
$ws.Range("K24").Value = 0.4986945169712794
$ws.Range("L24").Value = 191
$ws.Range("M24").Value = 191
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 192
$ws.Range("J25").Value = "like"
$ws.Range("K25").Value = 0.4735294117647059
$ws.Range("L25").Value = 161
$ws.Range("M25").Value = 161
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 179
$ws.Range("J26").Value = "care"
$ws.Range("K26").Value = 0.449438202247191
$ws.Range("L26").Value = 40
$ws.Range("M26").Value = 40
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 49
$ws.Range("J27").Value = "help"
$ws.Range("K27").Value = 0.3966101694915254
$ws.Range("L27").Value = 117
$ws.Range("M27").Value = 117
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 178
$ws.Range("J28").Value = "please"
$ws.Range("K28").Value = 0.3347280334728033
$ws.Range("L28").Value = 80
$ws.Range("M28").Value = 80
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 159
$ws.Range("J29").Value = "increase"
$ws.Range("K29").Value = 0.3205128205128205
$ws.Range("L29").Value = 25
$ws.Range("M29").Value = 25
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 53
$ws.Range("J30").Value = "you"
$ws.Range("K30").Value = 0.02333333333333333
$ws.Range("L30").Value = 28
$ws.Range("M30").Value = 28
$ws.Range("N30").Value = 1
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = $false
$ws.Range("Q30").Value = 1172
$ws.Range("J31").Value = "."
$ws.Range("K31").Value = 0.005204163330664532
$ws.Range("L31").Value = 26
$ws.Range("M31").Value = 26
$ws.Range("N31").Value = 1
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = $false
$ws.Range("Q31").Value = 4970
